$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (header + 17 player rows), replacing the previous
# 18-row data set. "Amen Thompson" was removed and the remaining rows
# were re-ordered.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Caris LeVert", "SG,SF", "Cleveland Cavaliers"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers")
)

$rowCount = $data.Length

# Clear out the old range first (old data went down to row 19) then
# write the new values.
$ws.Range("A1:C19").ClearContents()

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
